# Release-Notes.xlsx update
# - "Folder Inventory": a brand-new folder ("Get data into Fabric Lakehouse")
#   was discovered and is inserted as the new, most-recent row (row 7),
#   pushing every row that used to be at 7..70 down to 8..71.
# - "Metadata": refresh timestamp, bump Total Folders (69 -> 70) and
#   Workflow Run (text "4" -> text "5").
# - "Summary": bump Total Folders and Folders with Files (69 -> 70 each).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Folder Inventory: insert the new folder row at the top of the data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Folder Inventory")

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Get data into Fabric Lakehouse"
$ws.Range("B7").Value = "Get data into Fabric Lakehouse"
$ws.Range("C7").Value = "2025-06-11 15:00:50 +0000"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Root"

# ---------------------------------------------------------------------
# 2) Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2025-06-11 15:42:53 UTC"
$meta.Range("B4").Value = 70

# "Workflow Run" is stored as text ("4" -> "5"), not a number - force the
# cell to text formatting before writing so it round-trips as a string.
$runCell = $meta.Range("B5")
$runCell.NumberFormat = "@"
$runCell.Value = "5"

# ---------------------------------------------------------------------
# 3) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 70
$summary.Range("B3").Value = 70
